# Applies the scheduled-runner price/profit refresh to the Famfrit_Profits
# workbook. Each worksheet corresponds to a crafting class tab; rows are
# addressed by their existing row number (no rows are inserted or removed),
# only specific H/I/J/K/L/M/N cell values are refreshed, and a couple of
# cells are added/removed where the refreshed data no longer has a value.

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        $Sheet,
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $Values.Keys) {
        $addr = "$col$Row"
        $Sheet.Range($addr).Value = $Values[$col]
    }
}

function Clear-Cells {
    param(
        $Sheet,
        [int]$Row,
        [string[]]$Cols
    )
    foreach ($col in $Cols) {
        $addr = "$col$Row"
        $Sheet.Range($addr).ClearContents()
    }
}

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

Set-Cells $ws 4 @{ H=165.85715; I=165.85715; K=165.85715; M=-51.85714999999999 }

Set-Cells $ws 17 @{ H=989364.8; J=1004126; L=3012378; N=-3012714 }

Set-Cells $ws 40 @{ H=3062; J=3995; L=3995; N=-4345 }

Set-Cells $ws 47 @{ H=13249.5; I=19999; J=6500; K=19999; L=6500; M=-19027; N=-8444 }

Set-Cells $ws 100 @{ H=2091.8235; I=1822.6666; J=2394.625; K=1822.6666; L=2394.625; M=-1281.6666; N=-3476.625 }

Set-Cells $ws 110 @{ H=90000; J=90000; L=90000; N=-98180 }

Set-Cells $ws 127 @{ H=3357.2; I=650; J=4034; K=1950; L=12102; M=3010; N=-22022 }

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

Set-Cells $ws 2 @{ H=6124.091; I=5928.8887; K=5928.8887; M=-5815.8887 }

Set-Cells $ws 30 @{ H=0; J=0; L=0 }
Clear-Cells $ws 30 @("N")

Set-Cells $ws 45 @{ H=3950.111; I=3480.0833; J=4890.1665; K=3480.0833; L=4890.1665; M=-3103.0833; N=-5644.1665 }

Set-Cells $ws 110 @{ H=32766.545; I=35978.2; K=35978.2; M=-33933.2 }

Set-Cells $ws 116 @{ H=6124.091; I=5928.8887; K=5928.8887; M=-3634.8887 }

Set-Cells $ws 132 @{ H=69128.53999999999; I=4715.6; J=345184; K=14146.8; L=1035552; M=-11616.8; N=-1040612 }

Set-Cells $ws 135 @{ H=118575; J=118575; L=118575; N=-128715 }

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

Set-Cells $ws 3 @{ H=6124.091; I=5928.8887; K=5928.8887; M=-5814.8887 }

Set-Cells $ws 113 @{ H=3973.75; I=3973.75; K=3973.75; M=-1803.75 }

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-Cells $ws 5 @{ H=2496.1; I=492.4; J=4499.8; K=1477.2; L=13499.4; M=-1365.2; N=-13723.4 }

Set-Cells $ws 9 @{ H=7733.1816; I=75075; J=999; K=225225; L=2997; M=-225001; N=-3445 }

Set-Cells $ws 22 @{ H=2322.25; I=1732.25; J=2912.25; K=5196.75; L=8736.75; M=-5027.75; N=-9074.75 }

Set-Cells $ws 27 @{ H=2322.25; I=1732.25; J=2912.25; K=5196.75; L=8736.75; M=-5094.75; N=-8940.75 }

Set-Cells $ws 135 @{ H=2496.1; I=492.4; J=4499.8; K=4431.599999999999; L=40498.2; M=-1896.599999999999; N=-45568.2 }

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

Set-Cells $ws 16 @{ H=969.5; I=925.37036; K=925.37036; M=-755.37036 }

Set-Cells $ws 46 @{ H=1686.1875; J=2994.75; L=2994.75; N=-3370.75 }

Set-Cells $ws 61 @{ H=25002898; I=33335968; J=3689.2; K=33335968; L=3689.2; M=-33335766; N=-4093.2 }

Set-Cells $ws 93 @{ H=4299.3335; I=3999; K=3999; M=-2751 }

Set-Cells $ws 113 @{ H=25002898; I=33335968; J=3689.2; K=33335968; L=3689.2; M=-33333798; N=-8029.2 }

Set-Cells $ws 133 @{ H=55215; J=55215; L=55215; N=-60275 }

Set-Cells $ws 136 @{ H=2971.5; I=2812.3684; K=8437.1052; M=-5887.1052 }

Set-Cells $ws 140 @{ H=95999.39999999999; J=97499.5; L=97499.5; N=-107859.5 }

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")

Set-Cells $ws 113 @{ H=1085.5; I=1010.2222; K=3030.6666; M=-860.6666 }

Set-Cells $ws 137 @{ H=120045; I=0; J=120045; K=0; L=120045; N=-130245 }
Clear-Cells $ws 137 @("M")
